# Fixed a bug in FastPlugin
# The data rows (A2:F21) were shuffled back into the correct order.
# This reproduces the fix by writing the corrected values for each
# row of the reel-weighting table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(1001, 18, 30, 75, 60, 72),
    @(601,  9,  60, 67, 60, 42),
    @(201,  9,  30, 15, 45, 30),
    @(1203, 3,  15, 15, 15, 15),
    @(901,  16, 15, 45, 60, 60),
    @(301,  6,  45, 30, 60, 45),
    @(801,  3,  67, 65, 52, 45),
    @(401,  9,  48, 67, 75, 45),
    @(101,  9,  30, 15, 60, 15),
    @(902,  1,  0,  0,  0,  0),
    @(501,  9,  52, 30, 75, 45),
    @(701,  3,  90, 45, 97, 15),
    @(1201, 2,  10, 10, 10, 10),
    @(1202, 2,  10, 10, 10, 10),
    @(802,  0,  4,  5,  4,  0),
    @(1101, 0,  15, 30, 30, 0),
    @(2,    0,  2,  2,  2,  2),
    @(502,  0,  4,  0,  0,  0),
    @(1,    0,  2,  2,  2,  2),
    @(3,    0,  3,  3,  3,  3)
)

$startRow = 2
for ($i = 0; $i -lt $values.Length; $i++) {
    $r = $startRow + $i
    $row = $values[$i]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
